$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Brn_Code column (C) changes from text "195" to numeric Brn_Code values
$ws.Range("C2").Value = 193
$ws.Range("C3").Value = 19
$ws.Range("C4").Value = 19

# Update the active selection cell shown in the sheet view
$null = $ws.Range("F10").Select()
